$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 0.329814376738541
$ws.Range("B2").Value = 0.2954230338842676
$ws.Range("C2").Value = 0.3628265025940295
$ws.Range("D2").Value = 0.2721144324105971
$ws.Range("E2").Value = 0.4068229258834966
$ws.Range("K2").Value = 2.308700637169787
$ws.Range("L2").Value = 2.067961237189873
$ws.Range("M2").Value = 2.539785518158206
$ws.Range("N2").Value = 1.90480102687418
$ws.Range("O2").Value = 2.847760481184476
$ws.Range("P2").Value = 64.54441
$ws.Range("Q2").Value = 46.29748610134416
$ws.Range("R2").Value = 87.42864091694132
$ws.Range("S2").Value = 53.61109968745336
$ws.Range("T2").Value = 74.66188525398353
$ws.Range("U2").Value = 0.2028046299179496
$ws.Range("V2").Value = 0.164399412302132
$ws.Range("W2").Value = 0.2415665995288183
$ws.Range("X2").Value = 0.1359348926839613
$ws.Range("Y2").Value = 0.2797384063660836
$ws.Range("Z2").Value = 0.8569092516457457
$ws.Range("AA2").Value = 0.8109437027015787
$ws.Range("AB2").Value = 0.8977337412525318
$ws.Range("AC2").Value = 0.7705067053523447
$ws.Range("AD2").Value = 0.9300398675239179
$ws.Range("F3").Value = 7.001189101731514
$ws.Range("G3").Value = 6.212024365265301
$ws.Range("H3").Value = 7.754818811715914
$ws.Range("I3").Value = 5.703289538520241
$ws.Range("J3").Value = 8.785199554711253
$ws.Range("K3").Value = 2.3103924035714
$ws.Range("L3").Value = 2.04996804053755
$ws.Range("M3").Value = 2.559090207866252
$ws.Range("N3").Value = 1.882085547711679
$ws.Range("O3").Value = 2.899115853054714
$ws.Range("P3").Value = 63.52424
$ws.Range("Q3").Value = 56.39283059166551
$ws.Range("R3").Value = 72.2141371691701
$ws.Range("S3").Value = 59.37693894815393
$ws.Range("T3").Value = 67.81250866268077
$ws.Range("U3").Value = 0.2028134555222862
$ws.Range("V3").Value = 0.1615561728265415
$ws.Range("W3").Value = 0.2447438339307111
$ws.Range("X3").Value = 0.131546783205473
$ws.Range("Y3").Value = 0.286157498810601
$ws.Range("Z3").Value = 0.8560805618415306
$ws.Range("AA3").Value = 0.8064146390035332
$ws.Range("AB3").Value = 0.9001245295524177
$ws.Range("AC3").Value = 0.7632224970590397
$ws.Range("AD3").Value = 0.9341686915997631
$ws.Range("A4").Value = 0.3301008880891508
$ws.Range("B4").Value = 0.295501325460806
$ws.Range("C4").Value = 0.3632672235549657
$ws.Range("D4").Value = 0.2721877960524035
$ws.Range("E4").Value = 0.4070528584374838
$ws.Range("F4").Value = 7.001529098039073
$ws.Range("G4").Value = 6.219226625207799
$ws.Range("H4").Value = 7.753913499864607
$ws.Range("I4").Value = 5.709017536113828
$ws.Range("J4").Value = 8.768311080594755
$ws.Range("K4").Value = 2.311371204970385
$ws.Range("L4").Value = 1.959791563416804
$ws.Range("M4").Value = 2.650350248252937
$ws.Range("N4").Value = 1.721624650984149
$ws.Range("O4").Value = 3.117795841065266
$ws.Range("P4").Value = 65.34237
$ws.Range("Q4").Value = 45.59072657740801
$ws.Range("R4").Value = 94.65279409933459
$ws.Range("S4").Value = 53.69611793254276
$ws.Range("T4").Value = 76.44578313582664
$ws.Range("U4").Value = 0.2014640536348754
$ws.Range("V4").Value = 0.1451732985503876
$ws.Range("W4").Value = 0.2578642589046234
$ws.Range("X4").Value = 0.1030804296677801
$ws.Range("Y4").Value = 0.3126961887145306
$ws.Range("Z4").Value = 0.8490192072278532
$ws.Range("AA4").Value = 0.7773669948036431
$ws.Range("AB4").Value = 0.9082213497167891
$ws.Range("AC4").Value = 0.7058395560397965
$ws.Range("AD4").Value = 0.9520788124809904
